$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay text (match source inlineStr type)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated price / volume figures
$ws.Range("D2").Value = "27.317.85"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "1.775.81"
$ws.Range("E3").Value = "  +3.46%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "313.69"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "0.5144"
$ws.Range("E7").Value = "  +7.43%  "
$ws.Range("D8").Value = "0.3683"
$ws.Range("E8").Value = "  +6.46%  "
$ws.Range("D9").Value = "42.65"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").Value = "0.07398"
$ws.Range("E10").Value = "  +1.40%  "
$ws.Range("D11").Value = "1.088"
$ws.Range("E11").Value = "  +3.76%  "
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").Value = "20.49"
$ws.Range("E13").Value = "  +2.72%  "
$ws.Range("D14").Value = "6.069"
$ws.Range("E14").Value = "  +3.30%  "
$ws.Range("D15").Value = "1.768.85"
$ws.Range("E15").Value = "  +3.27%  "
$ws.Range("D16").Value = "6.958"
$ws.Range("D17").Value = "89.31"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").Value = "0.00001047"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").Value = "0.06436"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("D21").Value = "16.76"
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("D22").Value = "5.826"
$ws.Range("E22").Value = "  +3.10%  "
$ws.Range("D23").Value = "27.349.61"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").Value = "11.25"
$ws.Range("E24").Value = "  +3.71%  "
$ws.Range("D25").Value = "2.118"
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("D26").Value = "154.39"
$ws.Range("E26").Value = "  +1.56%  "
$ws.Range("D27").Value = "20.20"
$ws.Range("E27").Value = "  +2.44%  "
$ws.Range("D28").Value = "2.330"
$ws.Range("E28").Value = "  +10.81%  "
$ws.Range("D29").Value = "1.972.96"
$ws.Range("E29").Value = "  +3.34%  "
$ws.Range("D30").Value = "121.21"
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("D31").Value = "1.066"
$ws.Range("E31").Value = "  +4.05%  "
$ws.Range("D32").Value = "0.09764"
$ws.Range("E32").Value = "  +5.67%  "
$ws.Range("D33").Value = "5.580"
$ws.Range("E33").Value = "  +4.23%  "
$ws.Range("D34").Value = "3.627"
$ws.Range("E34").Value = "  +1.22%  "
$ws.Range("D35").Value = "0.02246"
$ws.Range("E35").Value = "  +1.87%  "
$ws.Range("D36").Value = "0.05970"
$ws.Range("E36").Value = "  +0.58%  "
$ws.Range("D37").Value = "11.26"
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("D38").Value = "0.6150"
$ws.Range("E38").Value = "  +3.26%  "
$ws.Range("D39").Value = "4.839"
$ws.Range("E39").Value = "  +1.34%  "
$ws.Range("D40").Value = "0.2023"
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("E41").Value = "  +0.98%  "
$ws.Range("D42").Value = "8.096"
$ws.Range("E42").Value = "  +7.92%  "
$ws.Range("D43").Value = "1.134"
$ws.Range("E43").Value = "  +3.15%  "
$ws.Range("D44").Value = "13.08"
$ws.Range("E44").Value = "  +2.76%  "
$ws.Range("D45").Value = "0.5772"
$ws.Range("E45").Value = "  +2.42%  "
$ws.Range("D46").Value = "3.634"
$ws.Range("E46").Value = "  +1.22%  "
$ws.Range("D47").Value = "121.59"
$ws.Range("E47").Value = "  +2.40%  "
$ws.Range("D48").Value = "1.889"
$ws.Range("E48").Value = "  +2.22%  "
$ws.Range("E49").Value = "  +2.67%  "
$ws.Range("D50").Value = "0.06715"
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("D51").Value = "70.65"
$ws.Range("E51").Value = "  +1.16%  "
